# Initial Commit Week 4
#
# A new Week-6 entry ("Scraping and visualising housing prices for
# different post codes in Vienna") is inserted at the top of the
# currently-empty portion of the plan, which pushes the existing
# Week 6-15 entries (rows 7-16) down one row each (to rows 8-17),
# carrying their row heights and the one hyperlink in column C along
# with them. Separately, a new entry is added for the previously-empty
# row 52 ("Lookback on the Gym Year (Gym Tracking)").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Shift the Title column (and row heights) for rows 7-17 down by one,
#     with the brand new entry landing in row 7 -------------------------
$newTitles = @{
    7  = "Scraping and visualising housing prices for different post codes in Vienna"
    8  = "Predicting the price of used cars"
    9  = "Loans (Data is Plural)"
    10 = "Procrastinated one about importance of retraining model (Day ahead spot prices)"
    11 = "Image classification: Cracks in concrete"
    12 = "Electricity something"
    13 = "Prisoners problem"
    14 = "Pytorch"
    15 = "Webscraping Digitec or Ebay Tech Prices (Focus on one, either smartphones, laptops or something else, maybe Apple)"
    16 = "F1 2022 vs 2021 data?"
    17 = "Luxembourg Stats (Compare weather with St. Gallen?)"
}

$newHeights = @{
    7  = 43.5
    8  = 14.5
    9  = 14.5
    10 = 43.5
    11 = 43.5
    12 = 14.5
    13 = 14.5
    14 = 14.5
    15 = 58
    16 = 14.5
    17 = 29
}

foreach ($r in 7..17) {
    $ws.Cells.Item($r, 2).Value2 = $newTitles[$r]
    $ws.Rows($r).RowHeight = $newHeights[$r]
}

# --- The hyperlink that used to sit in C10 (next to "Image classification:
#     Cracks in concrete") now belongs in C11, alongside its text, which
#     moved down with the rest of the block. This engine only supports
#     deleting the *entire* hyperlink collection of a sheet at once, so
#     clear it and rebuild all three links at their correct locations. ---
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Cells.Item(2, 3), "https://github.com/rfordatascience/tidytuesday/blob/master/data/2021/2021-07-27/readme.md") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5, 3), "https://www.kaggle.com/datasets/sinamhd9/concrete-comprehensive-strength?select=Concrete_Data.xls") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(11, 3), "https://www.kaggle.com/code/vishnu0399/ensuring-structural-safety-crack-detection") | Out-Null

# C10 no longer holds the hyperlinked text -> make sure it is blank now.
$ws.Cells.Item(10, 3).Value2 = $null

# --- Fill in the previously-empty row 52 with a new entry ---------------
$ws.Cells.Item(52, 2).Value2 = "Lookback on the Gym Year (Gym Tracking)"
$ws.Rows(52).RowHeight = 29

# --- View-state tweaks: zoom out a bit and move the selection to B53 ----
$ws.Application.ActiveWindow.Zoom = 145
$ws.Range("B53").Select()
